$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix header typos (Latin "B" -> Cyrillic "В") ---
$ws.Range("D1").Value = "SNR В0031-07 GURT"
$ws.Range("G1").Value = "SNR В0031-07 UTR2 sect 09"

# --- 2. Move the "DM B1133+16 UTR2 full" values from column O to column N ---
$moveRows = @(2, 3, 4, 5, 6, 8, 9, 10, 11, 12, 13)
foreach ($r in $moveRows) {
    $src = $ws.Cells.Item($r, 15)   # column O
    $dst = $ws.Cells.Item($r, 14)   # column N
    $dst.Value = $src.Value2
    $src.Clear()
}

# --- 3. Remove leftover sentinel/placeholder values ---
$ws.Range("V3").Clear()
$ws.Range("O7").Clear()
$ws.Range("U7").Clear()
$ws.Range("Q10").Clear()
$ws.Range("R10").Clear()
$ws.Range("S10").Clear()
$ws.Range("U10").Clear()
$ws.Range("AA10").Clear()
$ws.Range("AA11").Clear()
$ws.Range("S13").Clear()
$ws.Range("U13").Clear()
$ws.Range("V13").Clear()
$ws.Range("O14").Clear()
$ws.Range("S14").Clear()
$ws.Range("U14").Clear()
$ws.Range("V14").Clear()
$ws.Range("AA14").Clear()

# --- 4. Populate new plot data on row 24 ---
$ws.Range("E24").Value = 11.4
$ws.Range("F24").Value = 13.5
$ws.Range("H24").Value = 10.896
$ws.Range("L24").Value = 16.1
$ws.Range("R24").Value = 17.8
$ws.Range("X24").Value = 40

# Highlight W24 in yellow (no value change)
$ws.Range("W24").Interior.Color = 65535

# --- 5. Update the frozen pane / selection in the sheet view ---
$ws.Range("G1").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("G1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("N1").Select()
